# [Improvement] On terminology : room -> bed
$wb = $excel.ActiveWorkbook

# Rename the "rooms" worksheet to "beds"
$roomsSheet = $wb.Worksheets.Item("rooms")
$roomsSheet.Name = "beds"

# Update header row terminology on the renamed sheet
$roomsSheet.Range("A1").Value = "all_beds"
$roomsSheet.Range("B1").Value = "new_beds"
$roomsSheet.Range("C1").Value = "old_beds"
$roomsSheet.Range("E1").Value = "new_beds_service"
$roomsSheet.Range("F1").Value = "old_beds_service"
$roomsSheet.Range("G1").Value = "beds_capacities"

# "babies" sheet keeps its C1 selection but is no longer the active tab
$babiesSheet = $wb.Worksheets.Item("babies")
$babiesSheet.Activate()
$babiesSheet.Range("C1").Select()

# Make the "beds" sheet the active tab with its own selection
$roomsSheet.Activate()
$roomsSheet.Range("E22").Select()
